$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: rename company "Dennys" -> "Trivago"
$ws.Range("A5").Value = "Trivago"

# Row 5: update Category GUID string to a new unique value
$ws.Range("D5").Value = '"65e4a51476cc4294f12ba119"'
